$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 22:48:24"
$ws.Range("I2").Value = "4.3 mm"
$ws.Range("E3").Value = "2026-02-19 22:48:26"
$ws.Range("G3").Value = "254 cm"
$ws.Range("I3").Value = "6.4 mm"
$ws.Range("E4").Value = "2026-02-19 22:48:29"
$ws.Range("J4").Value = "1010.4 hPa"
$ws.Range("E5").Value = "2026-02-19 22:48:31"
$ws.Range("I5").Value = "8.3 mm"
$ws.Range("E6").Value = "2026-02-19 22:48:33"
$ws.Range("J6").Value = "1010.5 hPa"
$ws.Range("O6").Value = "10.2 °C"
$ws.Range("E7").Value = "2026-02-19 22:48:36"
$ws.Range("J7").Value = "1011.5 hPa"
$ws.Range("E8").Value = "2026-02-19 22:48:38"
$ws.Range("J8").Value = "1011.2 hPa"
$ws.Range("E9").Value = "2026-02-19 22:48:41"
$ws.Range("O9").Value = "10.2 °C"
$ws.Range("E10").Value = "2026-02-19 22:48:43"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "72%"
$ws.Range("N10").Value = "3.2 °C 22:29 TU"
$ws.Range("O10").Value = "9.8 °C"
$ws.Range("E11").Value = "2026-02-19 22:48:44"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "61%"
$ws.Range("O11").Value = "6.0 °C"
$ws.Range("E12").Value = "2026-02-19 22:48:45"
$ws.Range("O12").Value = "11.0 °C"
$ws.Range("E13").Value = "2026-02-19 22:48:46"
$ws.Range("J13").Value = "1011.8 hPa"
$ws.Range("E14").Value = "2026-02-19 22:48:47"
$ws.Range("E15").Value = "2026-02-19 22:48:49"
$ws.Range("O15").Value = "9.8 °C"
$ws.Range("E16").Value = "2026-02-19 22:48:50"
$ws.Range("I16").Value = "11.8 mm"
$ws.Range("E17").Value = "2026-02-19 22:48:51"
$ws.Range("E18").Value = "2026-02-19 22:48:52"
$ws.Range("E19").Value = "2026-02-19 22:48:53"
$ws.Range("O19").Value = "5.1 °C"
$ws.Range("E20").Value = "2026-02-19 22:48:54"
$ws.Range("L20").Value = "98.3 km/h - 328º 22:22 TU"
$ws.Range("E21").Value = "2026-02-19 22:48:55"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "60%"
$ws.Range("J21").Value = "1011.9 hPa"
$ws.Range("E22").Value = "2026-02-19 22:48:58"
$ws.Range("I22").Value = "0.5 mm"
$ws.Range("K22").Value = "15.5 MJ/m2"
$ws.Range("O22").Value = "-7.6 °C"
$ws.Range("E23").Value = "2026-02-19 22:49:00"
$ws.Range("I23").Value = "12.0 mm"
$ws.Range("E24").Value = "2026-02-19 22:49:03"
$ws.Range("J24").Value = "1015.5 hPa"
$ws.Range("E25").Value = "2026-02-19 22:49:05"
$ws.Range("I25").Value = "7.7 mm"
$ws.Range("E26").Value = "2026-02-19 22:49:08"
$ws.Range("J26").Value = "1010.3 hPa"
$ws.Range("L26").Value = "95.0 km/h - 346º 22:27 TU"
$ws.Range("O26").Value = "3.2 °C"
$ws.Range("E27").Value = "2026-02-19 22:49:11"
$ws.Range("E28").Value = "2026-02-19 22:49:13"
$ws.Range("J28").Value = "1010.4 hPa"
$ws.Range("O28").Value = "9.0 °C"
$ws.Range("E29").Value = "2026-02-19 22:49:16"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "78%"
$ws.Range("N29").Value = "4.2 °C 22:17 TU"
$ws.Range("O29").Value = "10.1 °C"
$ws.Range("E30").Value = "2026-02-19 22:49:18"
$ws.Range("J30").Value = "1010.5 hPa"
$ws.Range("E31").Value = "2026-02-19 22:49:21"
$ws.Range("L31").Value = "117.4 km/h - 346º 22:27 TU"
$ws.Range("O31").Value = "11.9 °C"
$ws.Range("E32").Value = "2026-02-19 22:49:23"
$ws.Range("E33").Value = "2026-02-19 22:49:26"
$ws.Range("J33").Value = "1011.4 hPa"
$ws.Range("O33").Value = "3.8 °C"
$ws.Range("E34").Value = "2026-02-19 22:49:28"
$ws.Range("E35").Value = "2026-02-19 22:49:31"
$ws.Range("J35").Value = "1016.9 hPa"
$ws.Range("E36").Value = "2026-02-19 22:49:33"
$ws.Range("J36").Value = "1010.7 hPa"
$ws.Range("O36").Value = "11.8 °C"
$ws.Range("E37").Value = "2026-02-19 22:49:36"
$ws.Range("J37").Value = "1011.9 hPa"
$ws.Range("O37").Value = "5.7 °C"
$ws.Range("E38").Value = "2026-02-19 22:49:38"
$ws.Range("O38").Value = "11.8 °C"
$ws.Range("E39").Value = "2026-02-19 22:49:41"
$ws.Range("E40").Value = "2026-02-19 22:49:43"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "69%"
$ws.Range("J40").Value = "1013.1 hPa"
$ws.Range("O40").Value = "6.6 °C"
$ws.Range("E41").Value = "2026-02-19 22:49:46"
$ws.Range("J41").Value = "1013.4 hPa"
$ws.Range("O41").Value = "14.0 °C"
$ws.Range("E42").Value = "2026-02-19 22:49:48"
$ws.Range("N42").Value = "5.3 °C 22:24 TU"
$ws.Range("O42").Value = "10.8 °C"
$ws.Range("E43").Value = "2026-02-19 22:49:50"
$ws.Range("N43").Value = "5.2 °C 22:28 TU"
$ws.Range("O43").Value = "8.9 °C"
$ws.Range("E44").Value = "2026-02-19 22:49:53"
$ws.Range("I44").Value = "10.4 mm"
$ws.Range("E45").Value = "2026-02-19 22:49:55"
$ws.Range("I45").Value = "3.6 mm"
$ws.Range("J45").Value = "1016.2 hPa"
$ws.Range("E46").Value = "2026-02-19 22:49:58"
$ws.Range("J46").Value = "1016.3 hPa"
